# Apply weekly refresh of Fruta/Hortaliza data: shuffle the D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg)
# values across rows 2-11 to reflect the latest values from origin.ss.digital feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, M, N, O, P, S
$data = @{
    2  = @{ D = 44798; M = 80;  N = 21000; O = 22000; P = 21500; S = 1075 }
    3  = @{ D = 45092; M = 150; N = 24000; O = 25000; P = 24333; S = 1217 }
    4  = @{ D = 44792; M = 100; N = 21000; O = 22000; P = 21500; S = 1075 }
    5  = @{ D = 44533; M = 100; N = 16000; O = 17000; P = 16500; S = 825 }
    6  = @{ D = 44890; M = 80;  N = 20000; O = 23000; P = 22250; S = 1112 }
    7  = @{ D = 44320; M = 80;  N = 16000; O = 17000; P = 16500; S = 825 }
    8  = @{ D = 44708; M = 80;  N = 20000; O = 21000; P = 20500; S = 1025 }
    9  = @{ D = 44357; M = 100; N = 14000; O = 15000; P = 14500; S = 725 }
    10 = @{ D = 44893; M = 80;  N = 21000; O = 22000; P = 21625; S = 1081 }
    11 = @{ D = 44761; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
}
